$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Style-fix copies (copy format+value from a stable donor cell before overwriting values) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("I14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("I14").Copy($ws.Range("D29"))
$ws.Range("K14").Copy($ws.Range("E29"))
$ws.Range("I14").Copy($ws.Range("G29"))
$ws.Range("K14").Copy($ws.Range("H29"))
$ws.Range("I14").Copy($ws.Range("D30"))
$ws.Range("K14").Copy($ws.Range("E30"))
$ws.Range("I14").Copy($ws.Range("G30"))
$ws.Range("K14").Copy($ws.Range("H30"))
$ws.Range("C14").Copy($ws.Range("F33"))

# --- Final numeric value assignments ---
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -64.705882352941
$ws.Range("C15").Value = 3
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 30
$ws.Range("K15").Value = -25
$ws.Range("L15").Value = 30.434782608695
$ws.Range("M15").Value = 7.142857142857
$ws.Range("N15").Value = -6.25
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 12
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -73.333333333333
$ws.Range("I16").Value = 281
$ws.Range("J16").Value = 403
$ws.Range("K16").Value = -30.272952853598
$ws.Range("L16").Value = -14.067278287461
$ws.Range("M16").Value = -13.003095975232
$ws.Range("N16").Value = -76.69983416252
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -38.461538461538
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 48
$ws.Range("H17").Value = -39.583333333333
$ws.Range("I17").Value = 504
$ws.Range("J17").Value = 598
$ws.Range("K17").Value = -15.71906354515
$ws.Range("L17").Value = 15.068493150684
$ws.Range("M17").Value = 66.887417218543
$ws.Range("N17").Value = 27.272727272727
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = -30
$ws.Range("L18").Value = 6.870229007633
$ws.Range("M18").Value = -48.148148148148
$ws.Range("N18").Value = -92.303463441451
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -18.75
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = -21.212121212121
$ws.Range("I19").Value = 670
$ws.Range("J19").Value = 891
$ws.Range("K19").Value = -24.803591470258
$ws.Range("L19").Value = -14.866581956798
$ws.Range("M19").Value = 45.021645021645
$ws.Range("N19").Value = -49.737434358589
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -34.615384615384
$ws.Range("I20").Value = 227
$ws.Range("J20").Value = 290
$ws.Range("K20").Value = -21.724137931034
$ws.Range("L20").Value = -29.283489096573
$ws.Range("M20").Value = 2.714932126696
$ws.Range("N20").Value = -89.0760346487
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -36.170212765957
$ws.Range("F21").Value = 119
$ws.Range("G21").Value = 182
$ws.Range("H21").Value = -34.615384615384
$ws.Range("I21").Value = 1858
$ws.Range("J21").Value = 2426
$ws.Range("K21").Value = -23.413025556471
$ws.Range("L21").Value = -8.47290640394
$ws.Range("M21").Value = 15.475450590428
$ws.Range("N21").Value = -72.998110739718
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 42
$ws.Range("K22").Value = -43.243243243243
$ws.Range("L22").Value = -48.148148148148
$ws.Range("M22").Value = 50
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -53.333333333333
$ws.Range("F24").Value = 82
$ws.Range("G24").Value = 126
$ws.Range("H24").Value = -34.920634920634
$ws.Range("I24").Value = 1286
$ws.Range("J24").Value = 1840
$ws.Range("K24").Value = -30.108695652173
$ws.Range("L24").Value = -29.880043620501
$ws.Range("M24").Value = 23.298178331735
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -61.538461538461
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = -65.51724137931
$ws.Range("I25").Value = 422
$ws.Range("J25").Value = 946
$ws.Range("K25").Value = -55.391120507399
$ws.Range("L25").Value = -54.574811625403
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = 19.047619047619
$ws.Range("F26").Value = 76
$ws.Range("G26").Value = 89
$ws.Range("H26").Value = -14.606741573033
$ws.Range("I26").Value = 959
$ws.Range("J26").Value = 1214
$ws.Range("K26").Value = -21.004942339374
$ws.Range("L26").Value = 6.319290465631
$ws.Range("M26").Value = 11.901983663944
$ws.Range("C27").Value = 3
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = -31.666666666666
$ws.Range("L27").Value = -6.818181818181
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 140
$ws.Range("I28").Value = 112
$ws.Range("K28").Value = -14.503816793893
$ws.Range("L28").Value = -17.647058823529
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = -87.5
$ws.Range("N29").Value = -98.181818181818
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -83.333333333333
$ws.Range("N30").Value = -98.039215686274
